# Correct the confusion-matrix "Pred False" column for the Simulator sheet:
#  - Truth Maybe -> Pred False (D22): 14 -> 15
#  - Truth False -> Pred False (D23): 239 -> 238
# Everything downstream (row/column totals, diffs, unweighted & weighted
# per-class metrics, binary metrics) recalculates automatically from these
# two input cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Simulator")

$ws.Range("D22").Value = 15
$ws.Range("D23").Value = 238

# Leave the cursor where the author left it after making the edit.
$ws.Range("I26").Select()
